$wb = $excel.ActiveWorkbook

# --- Add the new "ticketData" worksheet at the end of the tab strip ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ticketData"

# --- Column widths (best-effort match of authored widths) ---
$ws.Columns.Item(1).ColumnWidth = 16.721354166666668
$ws.Columns.Item(2).ColumnWidth = 32.166666666666664
$ws.Columns.Item(3).ColumnWidth = 25.498697916666668
$ws.Columns.Item(5).ColumnWidth = 31.608072916666668

# --- Cell values. Order chosen to reproduce the shared-string table order seen in the target file ---
$ws.Range("D1").Value = "client"
$ws.Range("E1").Value = "description"
$ws.Range("C2").Value = "Lekshmi ticket"
$ws.Range("E2").Value = "Lekshmi ticket added client abc"
$ws.Range("C3").Value = "rfs"
$ws.Range("D3").Value = "dfgsd"
$ws.Range("E3").Value = "sdhdgj"
$ws.Range("C1").Value = "tittle"
$ws.Range("D2").Value = "ABC"

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin@admin.com"
$ws.Range("B2").Value = 12345678
$ws.Range("A3").Value = "test"
$ws.Range("B3").Value = 344

# --- Update selections on the other sheets that changed ---
$wsClient = $wb.Worksheets.Item("clientData")
[void]$wsClient.Range("C17").Select()

$wsClockout = $wb.Worksheets.Item("clockoutData")
[void]$wsClockout.Range("B18").Select()

$wsAnnouncement = $wb.Worksheets.Item("announcementData")
[void]$wsAnnouncement.Range("A2:B3").Select()

$wsEvent = $wb.Worksheets.Item("eventData")
[void]$wsEvent.Range("A1:B4").Select()

# --- Selection / active cell + active tab on the new sheet (done last so it ends up the active sheet) ---
[void]$ws.Range("D2").Select()
